$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - best_params text updates
$ws.Range("F2").Value = "{'max_depth': 20, 'n_estimators': 50}"
$ws.Range("G2").Value = "{'learning_rate': 0.1, 'max_depth': 3, 'n_estimators': 200}"
$ws.Range("K2").Value = "{'activation': 'leaky_relu', 'b_random_vec_range': [0, 10], 'lam': 2, 'n_layer': 16, 'n_nodes': 128, 'random_seed': 911, 'same_feature': True, 'w_random_vec_range': [-10, 10]}"

# Row 3 - rmse updates
$ws.Range("F3").Value = 0.08881932908997739
$ws.Range("G3").Value = 0.09121419654056975
$ws.Range("H3").Value = 0.1049163036361207
$ws.Range("K3").Value = 0.0357367781029535

# Row 4 - r2 updates
$ws.Range("F4").Value = 0.4997195810875087
$ws.Range("G4").Value = 0.4562974940407168
$ws.Range("H4").Value = 0.3108353015635962
$ws.Range("K4").Value = 0.9079935044898744

# Row 5 - mape updates
$ws.Range("F5").Value = 37.82304280483131
$ws.Range("G5").Value = 32.47377899778003
$ws.Range("H5").Value = 48.29575720014626
$ws.Range("K5").Value = 11.78003579801157
